$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsMappings = $wb.Worksheets.Item("mappings")

# ---------------------------------------------------------------------------
# 1. Update the "address-data" JSON content (data sheet, B2) to add
#    cardinality information to the structure + its children.
# ---------------------------------------------------------------------------
$wsData.Range("B2").Value = '{"ident":"address-data","cardinality":"ONE","children":[{"DataValue":{"ident":"street","structureType":"STRING","cardinality":"ONE"}},{"DataValue":{"ident":"number","structureType":"STRING","cardinality":"ONE"}},{"DataValue":{"ident":"postcode","structureType":"NUMBER","cardinality":"ONE"}},{"DataValue":{"ident":"town","structureType":"STRING","cardinality":"ONE"}}]}'

# ---------------------------------------------------------------------------
# 2. Give the "mappings" sheet B2 cell (the address-mapping resolver entry)
#    a bigger/ different font - 14pt Times, black - and grow the row to fit.
# ---------------------------------------------------------------------------
$mapCell = $wsMappings.Range("B2")
$mapCell.Font.Name = "Times"
$mapCell.Font.Size = 14
$mapCell.Font.Color = 0

$wsMappings.Rows.Item(2).RowHeight = 19

# ---------------------------------------------------------------------------
# 3. Update the remembered cell-selections on both sheets (cosmetic sheet
#    view state) while keeping "mappings" as the active tab.
# ---------------------------------------------------------------------------
$wsData.Activate()
$wsData.Range("B2").Select()

$wsMappings.Activate()
$wsMappings.Range("C6").Select()
